$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = "10x_visiumhd"
$ws.Range("C3").Select()
